$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has 3 data rows (id 47, 48, 49 at rows 2,3,4). The
# target state has 8 data rows (id 54,49,50,51,47,48,52,53 at rows 2..9).
# Rows that keep exactly the same field values as before (id 47, 48, 49) are
# relocated using row copy/insert operations (which keeps their existing
# cell styling untouched). New rows are inserted blank and then populated
# with fresh values.
# ---------------------------------------------------------------------------

# Step 1: move id 49 (currently row 4) so it sits right after id 47 (row 2),
# i.e. make it the new row 3. Copy row4 and insert the copy before row3,
# then delete the now-duplicated original (shifted down to row5).
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(3).Insert(-4121)
$excel.CutCopyMode = 0
$ws.Rows.Item(5).Delete()

# Current order now: row2=47, row3=49, row4=48

# Step 2: make room for the two brand new rows (id 50, id 51) between
# row3 (id49) and row4 (id48).
$ws.Rows.Item(4).Insert(-4121)
$ws.Rows.Item(5).Insert(-4121)

# Current order now: row2=47, row3=49, row4=(blank/50), row5=(blank/51), row6=48

# Step 3: make room for two brand new rows (id 52, id 53) after id 48.
$ws.Rows.Item(7).Insert(-4121)
$ws.Rows.Item(8).Insert(-4121)

# Current order now: row2=47, row3=49, row4=(blank/50), row5=(blank/51),
#                     row6=48, row7=(blank/52), row8=(blank/53)

# Step 4: relocate id 47 (row2) so that it sits right before id 48, i.e.
# becomes the new row6 (pushing id48 and the two blank rows after it down by
# one). Copy row2 and insert the copy before row6's current position.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(6).Insert(-4121)
$excel.CutCopyMode = 0

# Current order now: row2=47(leftover original), row3=49, row4=(blank/50),
#                     row5=(blank/51), row6=47(new copy), row7=48,
#                     row8=(blank/52), row9=(blank/53)

# Row2's original content is now redundant (it was duplicated to row6), so
# row2 is free to be overwritten in place with the brand new id=54 record.

function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h, $j, $k, $l, $m, $n) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
}

# row2 : id 54
Set-Row 2 54 "Giuseppe" "Cangemi" "giuseppecangemi94@gmail.com" "derryrockpubfidelity.png" "SI" 7 "SI" "MOTTA SANT'ANASTASIA" 34679 "e" "Femminile" "eee"

# row4 : id 50
Set-Row 4 50 "alessandro " "aletta" "abc@prova.it" "qrcode_google_form.png" "SI" 3 "SI" "a" 34378 "a" "Maschile" "aaaa"

# row5 : id 51
Set-Row 5 51 "serena" "zante" "ss@ss.it" "derryrockpubfidelity.png" "SI" 4 "SI" "s" 34649 "e" "Maschile" "aaaa"

# row8 : id 52
Set-Row 8 52 "Giuseppe" "Cangemi" "giuseppecangemi94@gmail.com" "image.jpg" "SI" 5 "SI" "S" 34640 "Ss" "Maschile" "Der"

# row9 : id 53
Set-Row 9 53 "Giuseppe" "Cangemi" "giuseppecangemi94@gmail.com" "derryrockpubfidelity.png" "SI" 6 "SI" "MOTTA SANT'ANASTASIA" 34681 "e" "Femminile" "eee"

Write-Host "done"
